# Fixed typo in sample files: "License" -> "Licence" in the WMT_Extract
# header row, and restore the view (scroll position + selection) that was
# active when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WMT_Extract")

# Header row 1, columns W..AD held "LicenseTierX" - correct the spelling.
$ws.Range("W1").Value = "LicenceTier0"
$ws.Range("X1").Value = "LicenceTierD2"
$ws.Range("Y1").Value = "LicenceTierD1"
$ws.Range("Z1").Value = "LicenceTierC2"
$ws.Range("AA1").Value = "LicenceTierC1"
$ws.Range("AB1").Value = "LicenceTierB2"
$ws.Range("AC1").Value = "LicenceTierB1"
$ws.Range("AD1").Value = "LicenceTierA"

# Restore the sheet's scroll position and selection.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 21
$ws.Range("AD2").Select()
